$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "ETA" header in column E
$ws.Range("E1").Value = "ETA"

# Move the active selection to F1, matching the post-edit state
$ws.Range("F1").Select()
